# Add a new "override_normalization" column to the Trend_instructions sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in J1 and its value "T" in J2 (row 2 is the single data row).
$ws.Range("J1").Value = "override_normalization"
$ws.Range("J2").Value = "T"

# Move/keep the active selection on the newly added cell, matching the
# sheetView's selection recorded in the workbook.
$ws.Range("J2").Select()
